# Working on TMC429 implementation
# Update the "TMC429 Status" bitfield breakdown table (columns P:R, rows 3-10)
# on the "Tlm Struct" sheet: a new "On Target" flag row is introduced, the
# existing flags shift down by one, and a second TMC429 status-register
# example value (0x060504) is added alongside the existing 0x030201 one.
#
# Stable donor cells used purely to replicate exact cell-style indices
# (via copy/PasteSpecial-formats) without disturbing their own values:
#   C2  -> style used by Q-column text cells ("wrap text")
#   P2  -> style used by P-column number cells ("center/wrap")
#   C15 -> style used by most P/Q/R cells in this block ("center/wrap")
#   D15 -> style used by the register-value cells ("center/wrap")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tlm Struct")

# --- Row 3: RESERVED len 2 -> 1 ------------------------------------------
$ws.Range("P3").Value = 1

# --- Cells whose values only change to values/strings already used elsewhere
$ws.Range("Q5").Value = "Right limit"
$ws.Range("R5").Value = 1

$ws.Range("Q7").Value = "Homing Status"
$ws.Range("R7").Value = 3

$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = "TMC2130 Status"

$ws.Range("Q9").Value = "Target Pos"

$ws.Range("P10").Value = 24
$ws.Range("Q10").Value = "Actual Pos"
$ws.Range("C15").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = "0x030201"

# --- New shared string #1 (must be written before "On Target" below so the
#     two new shared-string table entries land in the expected order) ------
$ws.Range("D15").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = "0x060504"

# R8 becomes a plain number now (used to hold the "0x030201" string)
$ws.Range("C15").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 1

# --- New shared string #2 --------------------------------------------------
$ws.Range("Q4").Value = "On Target"
$ws.Range("C2").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 0

# --- Row 6: was "Homing Status" 2/3 -> becomes "Left limit" 1/0 -----------
$ws.Range("P2").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 1

$ws.Range("C2").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = "Left limit"
$ws.Range("R6").Value = 0

# --- Selection moves from R6 to S4 ----------------------------------------
$ws.Range("S4").Select()
